$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "REX_DEF" column (F) after the existing REX_DESC column (E).
$ws.Range("F1").Value = "REX_DEF"

# Match the formatting of the other header cells (bold, centered, bordered).
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Populate the new column's data rows with the REX_DEF values.
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
